$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1333.6666
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 1667.3334
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 1667.3334
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -2017.3334
$ws.Range("H137").Value = 4117.5713
$ws.Range("J137").Value = 6000.3335
$ws.Range("L137").Value = 18001.0005
$ws.Range("N137").Value = -23101.0005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 21595.4
$ws.Range("I45").Value = 50994
$ws.Range("J45").Value = 1996.3334
$ws.Range("K45").Value = 50994
$ws.Range("L45").Value = 1996.3334
$ws.Range("M45").Value = -50617
$ws.Range("N45").Value = -2750.3334
$ws.Range("H61").Value = 6287.5884
$ws.Range("I61").Value = 5430.5625
$ws.Range("K61").Value = 5430.5625
$ws.Range("M61").Value = -5218.5625
$ws.Range("H102").Value = 4305.591
$ws.Range("I102").Value = 2170.3125
$ws.Range("J102").Value = 9999.666999999999
$ws.Range("K102").Value = 2170.3125
$ws.Range("L102").Value = 9999.666999999999
$ws.Range("M102").Value = -548.3125
$ws.Range("N102").Value = -13243.667
$ws.Range("H110").Value = 862.8929000000001
$ws.Range("I110").Value = 789.0833
$ws.Range("J110").Value = 1305.75
$ws.Range("K110").Value = 789.0833
$ws.Range("L110").Value = 1305.75
$ws.Range("M110").Value = 1255.9167
$ws.Range("N110").Value = -5395.75
$ws.Range("H122").Value = 2032.1177
$ws.Range("I122").Value = 1472.5
$ws.Range("K122").Value = 4417.5
$ws.Range("M122").Value = -1967.5
$ws.Range("H132").Value = 1564.9131
$ws.Range("I132").Value = 1564.9131
$ws.Range("K132").Value = 4694.7393
$ws.Range("M132").Value = -2164.7393
$ws.Range("H136").Value = 6287.5884
$ws.Range("I136").Value = 5430.5625
$ws.Range("K136").Value = 16291.6875
$ws.Range("M136").Value = -13741.6875

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1195.3334
$ws.Range("J22").Value = 126.5
$ws.Range("L22").Value = 126.5
$ws.Range("N22").Value = -472.5
$ws.Range("H105").Value = 4068.5227
$ws.Range("I105").Value = 4649.057
$ws.Range("J105").Value = 1810.8889
$ws.Range("K105").Value = 4649.057
$ws.Range("L105").Value = 1810.8889
$ws.Range("M105").Value = -2902.057
$ws.Range("N105").Value = -5304.8889
$ws.Range("H107").Value = 6221.222
$ws.Range("I107").Value = 1497.75
$ws.Range("K107").Value = 1497.75
$ws.Range("M107").Value = 422.25

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1862.25
$ws.Range("I58").Value = 1634.8
$ws.Range("K58").Value = 1634.8
$ws.Range("M58").Value = -1431.8
$ws.Range("H69").Value = 9905.444
$ws.Range("I69").Value = 5733.1665
$ws.Range("J69").Value = 18250
$ws.Range("K69").Value = 5733.1665
$ws.Range("L69").Value = 18250
$ws.Range("M69").Value = -4984.1665
$ws.Range("N69").Value = -19748
$ws.Range("H72").Value = 9905.444
$ws.Range("I72").Value = 5733.1665
$ws.Range("J72").Value = 18250
$ws.Range("K72").Value = 17199.4995
$ws.Range("L72").Value = 54750
$ws.Range("M72").Value = -13455.4995
$ws.Range("N72").Value = -62238
$ws.Range("H99").Value = 3388
$ws.Range("I99").Value = 2250
$ws.Range("J99").Value = 3957
$ws.Range("K99").Value = 2250
$ws.Range("L99").Value = 3957
$ws.Range("M99").Value = -752
$ws.Range("N99").Value = -6953
$ws.Range("H126").Value = 3388
$ws.Range("I126").Value = 2250
$ws.Range("J126").Value = 3957
$ws.Range("K126").Value = 6750
$ws.Range("L126").Value = 11871
$ws.Range("M126").Value = -4280
$ws.Range("N126").Value = -16811
$ws.Range("H136").Value = 1862.25
$ws.Range("I136").Value = 1634.8
$ws.Range("K136").Value = 4904.4
$ws.Range("M136").Value = -2354.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2845.0977
$ws.Range("J131").Value = 3237.9697
$ws.Range("L131").Value = 9713.909100000001
$ws.Range("N131").Value = -19793.9091

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 38250
$ws.Range("J34").Value = 38250
$ws.Range("L34").Value = 38250
$ws.Range("N34").Value = -38786
$ws.Range("H52").Value = 30000
$ws.Range("J52").Value = 30000
$ws.Range("L52").Value = 30000
$ws.Range("N52").Value = -30518
$ws.Range("H76").Value = 38250
$ws.Range("J76").Value = 38250
$ws.Range("L76").Value = 38250
$ws.Range("N76").Value = -38880
$ws.Range("H79").Value = 38250
$ws.Range("J79").Value = 38250
$ws.Range("L79").Value = 38250
$ws.Range("N79").Value = -40434
$ws.Range("H122").Value = 15310.619
$ws.Range("I122").Value = 17472.4
$ws.Range("J122").Value = 9906.166999999999
$ws.Range("K122").Value = 52417.2
$ws.Range("L122").Value = 29718.501
$ws.Range("M122").Value = -49967.2
$ws.Range("N122").Value = -34618.501
$ws.Range("H132").Value = 2664.7693
$ws.Range("I132").Value = 2720.5833
$ws.Range("J132").Value = 1995
$ws.Range("K132").Value = 8161.749899999999
$ws.Range("L132").Value = 5985
$ws.Range("M132").Value = -5631.749899999999
$ws.Range("N132").Value = -11045

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 926
$ws.Range("I16").Value = 806.6667
$ws.Range("K16").Value = 806.6667
$ws.Range("M16").Value = -636.6667
$ws.Range("H22").Value = 2863.7778
$ws.Range("I22").Value = 1955.6
$ws.Range("K22").Value = 1955.6
$ws.Range("M22").Value = -1660.6
$ws.Range("H27").Value = 2863.7778
$ws.Range("I27").Value = 1955.6
$ws.Range("K27").Value = 1955.6
$ws.Range("M27").Value = -1848.6
$ws.Range("H43").Value = 19998.75
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 19998.75
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 19998.75
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -20384.75
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 4753.5
$ws.Range("J15").Value = 4753.5
$ws.Range("L15").Value = 4753.5
$ws.Range("N15").Value = -5329.5
$ws.Range("H40").Value = 13500
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 25000
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 25000
$ws.Range("M40").Value = -1851
$ws.Range("N40").Value = -25298
$ws.Range("H96").Value = 6375
$ws.Range("J96").Value = 5833.3335
$ws.Range("L96").Value = 5833.3335
$ws.Range("N96").Value = -8579.333500000001
